$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: "Sort" rows 5 and 6 by the Cylinder column (A), keeping the
# grouped vs/am (B,C) values attached to their own row so the later merge
# still makes sense. Column A and D:G swap between the two rows; B and C
# (already identical within each row's group) stay put.
$a5 = $ws.Range("A5").Value2
$a6 = $ws.Range("A6").Value2
$ws.Range("A5").Value = $a6
$ws.Range("A6").Value = $a5

$dg5 = $ws.Range("D5:G5").Value2
$dg6 = $ws.Range("D6:G6").Value2
$ws.Range("D5:G5").Value = $dg6
$ws.Range("D6:G6").Value = $dg5

# --- Step 2: rows 6 & 7 now share the same vs/am group (B=1,C=1) and the
# same Cylinder value (6). Vertically top-align the soon-to-be-merged anchor
# cells first (while they are still single cells) ...
$ws.Range("A6").VerticalAlignment = -4160
$ws.Range("B6").VerticalAlignment = -4160
$ws.Range("C6").VerticalAlignment = -4160

# ... then merge those columns vertically across the two rows.
$ws.Range("A6:A7").Merge()
$ws.Range("B6:B7").Merge()
$ws.Range("C6:C7").Merge()

# The un-merged lower-left cell (A7) should fall back to the plain "no data"
# look used elsewhere in the sheet (e.g. column H) rather than keep the
# heavier bordered/filled look it had before the merge.
$ws.Range("H7").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").ClearContents()

# B7/C7 are now silent, merged-away continuation cells -- drop their
# leftover value and formatting entirely.
$ws.Range("B7").ClearContents()
$ws.Range("B7").ClearFormats()
$ws.Range("C7").ClearContents()
$ws.Range("C7").ClearFormats()
